# Refresh the cryptocurrency price/volume snapshot in columns D (Price)
# and E (Volume(1h)) for the rows whose scraped values moved.
#
# Note: several "Price" strings (e.g. "1.003") look like plain decimals,
# so a leading apostrophe is used to keep them stored as text (matching
# the original inlineStr cells) instead of being auto-converted to a
# numeric value by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.119.17"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.668.65"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D5").Value = "'210.39"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").Value = "'0.5209"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").Value = "'0.2638"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "'0.06222"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").Value = "'0.07508"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "1.736.05"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "'4.422"
$ws.Range("E13").Value = "  -2.70%  "
$ws.Range("D14").Value = "'0.5587"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "'66.28"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "'0.000007892"
$ws.Range("E16").Value = "  -6.04%  "
$ws.Range("D17").Value = "26.162.70"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("D19").Value = "'4.779"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").Value = "'186.69"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("E21").Value = "  -5.69%  "
$ws.Range("D22").Value = "'6.173"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "'147.43"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'0.1243"
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("D26").Value = "'7.570"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("D27").Value = "'15.94"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'0.06215"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("D31").Value = "'3.477"
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").Value = "'3.422"
$ws.Range("E32").Value = "  -5.17%  "
$ws.Range("D33").Value = "'1.619"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("D34").Value = "'0.9939"
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("D35").Value = "'0.6034"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "'2.403"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'2.703"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'6.120"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").Value = "1.073.72"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'99.11"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "1.817.67"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "'55.94"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'0.05253"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "'7.917"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("D51").Value = "'5.938"
$ws.Range("E51").Value = "  -2.78%  "
